# OpenEMRData.xlsx edit script
# 1) Add a 4th row of data to the existing "InvalidCredentialTest" sheet.
# 2) Add two new worksheets: "AddPatientTest" and "AboutUsHeaderAndVersionTest"
#    with their own header/data rows (object-array style test data),
#    matching the "excel to objectarray, about us section" commit.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: InvalidCredentialTest - append a new test-data row
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "InvalidCredentialTest"

$ws1.Range("A4").Value = "bal"
$ws1.Range("B4").Value = "bal123"
$ws1.Range("C4").Value = "Danish"
$ws1.Range("D4").Value = "Invalid username or password123"

$ws1.Range("A1:C1").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: AddPatientTest (new sheet, placed right after sheet 1)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "AddPatientTest"

$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"
$ws2.Range("C1").Value = "Language"
$ws2.Range("D1").Value = "FirstName"
$ws2.Range("E1").Value = "LastName"
$ws2.Range("F1").NumberFormat = "@"
$ws2.Range("F1").Value = "DOB"
$ws2.Range("G1").Value = "Gender"
$ws2.Range("H1").Value = "ExpectedAlertText"
$ws2.Range("I1").Value = "ExpectedPatientName"

$ws2.Range("A2").Value = "admin"
$ws2.Range("B2").Value = "pass"
$ws2.Range("C2").Value = "English (Indian)"
$ws2.Range("D2").Value = "John"
$ws2.Range("E2").Value = "Ken"
$ws2.Range("F2").NumberFormat = "@"
$ws2.Range("F2").Value = "2021-07-19"
$ws2.Range("G2").Value = "Male"
$ws2.Range("H2").Value = "Assessment: Tobacco"
$ws2.Range("I2").Value = "Medical Record Dashboard - John Ken"

$ws2.Range("A3").Value = "admin"
$ws2.Range("B3").Value = "pass"
$ws2.Range("C3").Value = "English (Indian)"
$ws2.Range("D3").Value = "Bala"
$ws2.Range("E3").Value = "Dina"
$ws2.Range("F3").NumberFormat = "@"
$ws2.Range("F3").Value = "2021-07-19"
$ws2.Range("G3").Value = "Male"
$ws2.Range("H3").Value = "Assessment: Tobacco"
$ws2.Range("I3").Value = "Medical Record Dashboard - Bala Dina"

$ws2.Columns.Item(1).ColumnWidth = 10
$ws2.Columns.Item(2).ColumnWidth = 9.42578125
$ws2.Columns.Item(3).ColumnWidth = 14.85546875
$ws2.Columns.Item(6).ColumnWidth = 10.42578125
$ws2.Columns.Item(8).ColumnWidth = 17.7109375
$ws2.Columns.Item(9).ColumnWidth = 34.85546875

$ws2.Range("A1:C2").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 3: AboutUsHeaderAndVersionTest (new sheet, placed after AddPatientTest)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Add($null, $ws2)
$ws3.Name = "AboutUsHeaderAndVersionTest"

$ws3.Range("A1").Value = "Username"
$ws3.Range("B1").Value = "Password"
$ws3.Range("C1").Value = "Language"
$ws3.Range("D1").Value = "ExpectedHeader"
$ws3.Range("E1").Value = "ExpectedVersion"

$ws3.Range("A2").Value = "admin"
$ws3.Range("B2").Value = "pass"
$ws3.Range("C2").Value = "English (Indian)"
$ws3.Range("D2").Value = "About OpenEMR"
$ws3.Range("E2").Value = "Version Number: v6.0.0 (1)"

$ws3.Range("A3").Value = "physician"
$ws3.Range("B3").Value = "physician"
$ws3.Range("C3").Value = "English (Indian)"
$ws3.Range("D3").Value = "About OpenEMR"
$ws3.Range("E3").Value = "Version Number: v6.0.0 (1)"

$ws3.Columns.Item(3).ColumnWidth = 14.85546875
$ws3.Columns.Item(4).ColumnWidth = 15.7109375
$ws3.Columns.Item(5).ColumnWidth = 16.140625

$ws3.Range("A2").Select() | Out-Null

# Make the new "About Us" sheet the active / selected tab, as in the target file.
$ws3.Select() | Out-Null
